$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-18: update Price (D) / Volume(1h) (E) values in place ---
# NumberFormat "@" is applied first on cells whose new value would
# otherwise be auto-parsed as a number, to keep them text like the
# original inlineStr cells.
$ws.Range("D2").Value = "26.180.53"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.589.16"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.70"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.01"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "1.813.27"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.604.15"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.62"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "26.168.06"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  -0.73%  "

# --- Rows 19 & 20: Chainlink and BitcoinCash swap places (rank unchanged) ---
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.67"
$ws.Range("E19").Value = "  +1.64%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("E20").Value = "  -1.83%  "

# --- Rows 21-51: update Price (D) / Volume(1h) (E) values in place ---
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.00"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.47"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.96"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.06"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.16"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "1.414.98"
$ws.Range("E33").Value = "  +7.94%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.87"
$ws.Range("E40").Value = "  +4.54%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.946"
$ws.Range("E42").Value = "  -13.47%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").Value = "1.724.85"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.01"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.39"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  -0.19%  "
